$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2: 0 -> 1
$ws.Range("A2").Value = 1

# Update A3: 2 -> 0, B3: 1 -> 2
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 2

# Remove row 4 entirely (A4/B4), shrinking the used range to A1:B3
$ws.Rows.Item(4).Delete()
